{"js": "// Update the worksheet date heading and all 100 arithmetic problems in the\n// 20x5 table with new values, preserving all existing run/paragraph\n// formatting (fonts, sizes, alignment, etc).\n\n// --- 1) Update the date paragraph (\"2024-05-19 Sunday\" -> \"2024-05-20 Monday\") ---\nconst paras = context.document.body.paragraphs;\nparas.load(\"items\");\nawait context.sync();\n\nconst datePara = paras.items[0];\ndatePara.load(\"text\");\nawait context.sync();\n\nif (datePara.text === \"2024-05-19 Sunday\") {\n  datePara.insertText(\"2024-05-20 Monday\", Word.InsertLocation.replace);\n} else {\n  // Fallback: search the whole document for the old date text.\n  const results = context.document.body.search(\"2024-05-19 Sunday\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(\"2024-05-20 Monday\", Word.InsertLocation.replace);\n  }\n}\nawait context.sync();\n\n// --- 2) Update the 20x5 table of addition/subtraction problems ---\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"values\");\nawait context.sync();\n\n// Old -> new values, listed in row-major (top-left to bottom-right) order,\n// exactly matching the table's current reading order.\nconst oldValues = [\n  [\"39+58=\", \"90-63=\", \"93-79=\", \"96-63=\", \"74-64=\"],\n  [\"16+39=\", \"60+14=\", \"57-57=\", \"5+36=\", \"75-19=\"],\n  [\"34+35=\", \"86+5=\", \"80-23=\", \"8+37=\", \"64-62=\"],\n  [\"92-1=\", \"64-61=\", \"62+11=\", \"9+19=\", \"18+25=\"],\n  [\"9+69=\", \"13-12=\", \"23-14=\", \"61-36=\", \"60+19=\"],\n  [\"9+21=\", \"17-6=\", \"98-45=\", \"60+5=\", \"99-61=\"],\n  [\"34+21=\", \"43+49=\", \"6+82=\", \"40-9=\", \"47+29=\"],\n  [\"17-17=\", \"98-10=\", \"69-49=\", \"10+24=\", \"64-12=\"],\n  [\"95-57=\", \"60+4=\", \"14+53=\", \"98-10=\", \"24+6=\"],\n  [\"7+19=\", \"36-30=\", \"94-40=\", \"76+20=\", \"4+36=\"],\n  [\"19+18=\", \"68+24=\", \"19+35=\", \"6+92=\", \"54-21=\"],\n  [\"49+13=\", \"36-34=\", \"80-25=\", \"61-45=\", \"63-21=\"],\n  [\"69-44=\", \"84+4=\", \"42+25=\", \"3+69=\", \"71-18=\"],\n  [\"3+21=\", \"20+12=\", \"63-28=\", \"80+12=\", \"42+48=\"],\n  [\"87-86=\", \"19+35=\", \"51-18=\", \"11+5=\", \"5+46=\"],\n  [\"97-7=\", \"53+3=\", \"27+3=\", \"93-92=\", \"11+21=\"],\n  [\"6+13=\", \"26-14=\", \"24+46=\", \"23+58=\", \"17+27=\"],\n  [\"94-77=\", \"28+10=\", \"13+5=\", \"62-43=\", \"92-57=\"],\n  [\"33-32=\", \"77+18=\", \"0+75=\", \"28+54=\", \"10+58=\"],\n  [\"30-3=\", \"42+53=\", \"35-25=\", \"92-82=\", \"48-41=\"]\n];\n\nconst newValues = [\n  [\"11+73=\", \"24+74=\", \"27+12=\", \"74-21=\", \"75-46=\"],\n  [\"96-13=\", \"20+51=\", \"24+68=\", \"46+53=\", \"77-27=\"],\n  [\"63-17=\", \"66-59=\", \"61+26=\", \"43+53=\", \"59+5=\"],\n  [\"34-23=\", \"31-26=\", \"62-55=\", \"13+51=\", \"73+26=\"],\n  [\"82+1=\", \"75-0=\", \"40-33=\", \"49+19=\", \"61-17=\"],\n  [\"56+42=\", \"24-9=\", \"49+39=\", \"91+2=\", \"40+26=\"],\n  [\"80+14=\", \"48+7=\", \"91-69=\", \"7+18=\", \"30+7=\"],\n  [\"70+0=\", \"86-13=\", \"24-1=\", \"58+34=\", \"68-8=\"],\n  [\"27+10=\", \"5+94=\", \"10+51=\", \"67+22=\", \"10+82=\"],\n  [\"19+20=\", \"28-23=\", \"79-45=\", \"22+14=\", \"59-33=\"],\n  [\"44+22=\", \"47+30=\", \"7+86=\", \"24-10=\", \"94-23=\"],\n  [\"13+13=\", \"43-36=\", \"16+75=\", \"75-69=\", \"39+57=\"],\n  [\"65+32=\", \"83-68=\", \"89-20=\", \"58-47=\", \"1+61=\"],\n  [\"45-23=\", \"97-48=\", \"28+2=\", \"80-50=\", \"37+56=\"],\n  [\"36+49=\", \"41+24=\", \"12+29=\", \"48+5=\", \"11+48=\"],\n  [\"32+21=\", \"31-5=\", \"89-64=\", \"23+20=\", \"95-31=\"],\n  [\"66-29=\", \"18+38=\", \"77-75=\", \"49-13=\", \"55-4=\"],\n  [\"52-14=\", \"61-13=\", \"8+79=\", \"79-65=\", \"34+15=\"],\n  [\"33+33=\", \"62-28=\", \"26+69=\", \"96-17=\", \"75+12=\"],\n  [\"95-65=\", \"34-16=\", \"12+15=\", \"91-54=\", \"54-23=\"]\n];\n\nconst currentValues = table.values;\nconst updatedValues = currentValues.map((row, r) =>\n  row.map((cellText, c) => {\n    // Trust position: if the current cell still holds the expected old\n    // value use the mapped replacement, otherwise leave untouched.\n    if (oldValues[r] && oldValues[r][c] !== undefined && cellText === oldValues[r][c]) {\n      return newValues[r][c];\n    }\n    return cellText;\n  })\n);\n\ntable.values = updatedValues;\nawait context.sync();\n", "ps1": "# Update the worksheet date heading and all 100 arithmetic problems in the\n# 20x5 table with new values, preserving all existing run/paragraph\n# formatting (fonts, sizes, alignment, etc).\n\n$d = $word.ActiveDocument\n\n# --- 1) Update the date paragraph (\"2024-05-19 Sunday\" -> \"2024-05-20 Monday\") ---\n$dateOld = \"2024-05-19 Sunday\"\n$dateNew = \"2024-05-20 Monday\"\n\n$p1 = $d.Paragraphs.Item(1)\n$r1 = $p1.Range\nif ($r1.Text -eq ($dateOld + [char]13)) {\n    $r1.Text = $dateNew\n} elseif ($r1.Text -eq $dateOld) {\n    $r1.Text = $dateNew\n} else {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $dateOld\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $dateNew\n    $find.Execute([ref]$dateOld, $false, $false, $false, $false, $false, $true, 1, $false, $dateNew, 2)\n}\n\n# --- 2) Update the 20x5 table of addition/subtraction problems ---\n# Old -> new values, listed in row-major (top-left to bottom-right) order,\n# exactly matching the table's current reading order.\n$oldValues = @(\n    @(\"39+58=\", \"90-63=\", \"93-79=\", \"96-63=\", \"74-64=\"),\n    @(\"16+39=\", \"60+14=\", \"57-57=\", \"5+36=\", \"75-19=\"),\n    @(\"34+35=\", \"86+5=\", \"80-23=\", \"8+37=\", \"64-62=\"),\n    @(\"92-1=\", \"64-61=\", \"62+11=\", \"9+19=\", \"18+25=\"),\n    @(\"9+69=\", \"13-12=\", \"23-14=\", \"61-36=\", \"60+19=\"),\n    @(\"9+21=\", \"17-6=\", \"98-45=\", \"60+5=\", \"99-61=\"),\n    @(\"34+21=\", \"43+49=\", \"6+82=\", \"40-9=\", \"47+29=\"),\n    @(\"17-17=\", \"98-10=\", \"69-49=\", \"10+24=\", \"64-12=\"),\n    @(\"95-57=\", \"60+4=\", \"14+53=\", \"98-10=\", \"24+6=\"),\n    @(\"7+19=\", \"36-30=\", \"94-40=\", \"76+20=\", \"4+36=\"),\n    @(\"19+18=\", \"68+24=\", \"19+35=\", \"6+92=\", \"54-21=\"),\n    @(\"49+13=\", \"36-34=\", \"80-25=\", \"61-45=\", \"63-21=\"),\n    @(\"69-44=\", \"84+4=\", \"42+25=\", \"3+69=\", \"71-18=\"),\n    @(\"3+21=\", \"20+12=\", \"63-28=\", \"80+12=\", \"42+48=\"),\n    @(\"87-86=\", \"19+35=\", \"51-18=\", \"11+5=\", \"5+46=\"),\n    @(\"97-7=\", \"53+3=\", \"27+3=\", \"93-92=\", \"11+21=\"),\n    @(\"6+13=\", \"26-14=\", \"24+46=\", \"23+58=\", \"17+27=\"),\n    @(\"94-77=\", \"28+10=\", \"13+5=\", \"62-43=\", \"92-57=\"),\n    @(\"33-32=\", \"77+18=\", \"0+75=\", \"28+54=\", \"10+58=\"),\n    @(\"30-3=\", \"42+53=\", \"35-25=\", \"92-82=\", \"48-41=\")\n)\n\n$newValues = @(\n    @(\"11+73=\", \"24+74=\", \"27+12=\", \"74-21=\", \"75-46=\"),\n    @(\"96-13=\", \"20+51=\", \"24+68=\", \"46+53=\", \"77-27=\"),\n    @(\"63-17=\", \"66-59=\", \"61+26=\", \"43+53=\", \"59+5=\"),\n    @(\"34-23=\", \"31-26=\", \"62-55=\", \"13+51=\", \"73+26=\"),\n    @(\"82+1=\", \"75-0=\", \"40-33=\", \"49+19=\", \"61-17=\"),\n    @(\"56+42=\", \"24-9=\", \"49+39=\", \"91+2=\", \"40+26=\"),\n    @(\"80+14=\", \"48+7=\", \"91-69=\", \"7+18=\", \"30+7=\"),\n    @(\"70+0=\", \"86-13=\", \"24-1=\", \"58+34=\", \"68-8=\"),\n    @(\"27+10=\", \"5+94=\", \"10+51=\", \"67+22=\", \"10+82=\"),\n    @(\"19+20=\", \"28-23=\", \"79-45=\", \"22+14=\", \"59-33=\"),\n    @(\"44+22=\", \"47+30=\", \"7+86=\", \"24-10=\", \"94-23=\"),\n    @(\"13+13=\", \"43-36=\", \"16+75=\", \"75-69=\", \"39+57=\"),\n    @(\"65+32=\", \"83-68=\", \"89-20=\", \"58-47=\", \"1+61=\"),\n    @(\"45-23=\", \"97-48=\", \"28+2=\", \"80-50=\", \"37+56=\"),\n    @(\"36+49=\", \"41+24=\", \"12+29=\", \"48+5=\", \"11+48=\"),\n    @(\"32+21=\", \"31-5=\", \"89-64=\", \"23+20=\", \"95-31=\"),\n    @(\"66-29=\", \"18+38=\", \"77-75=\", \"49-13=\", \"55-4=\"),\n    @(\"52-14=\", \"61-13=\", \"8+79=\", \"79-65=\", \"34+15=\"),\n    @(\"33+33=\", \"62-28=\", \"26+69=\", \"96-17=\", \"75+12=\"),\n    @(\"95-65=\", \"34-16=\", \"12+15=\", \"91-54=\", \"54-23=\")\n)\n\n$t = $d.Tables.Item(1)\nfor ($r = 1; $r -le 20; $r++) {\n    for ($c = 1; $c -le 5; $c++) {\n        $cell = $t.Cell($r, $c)\n        $cellRange = $cell.Range\n        $expectedOld = $oldValues[$r - 1][$c - 1]\n        $newVal = $newValues[$r - 1][$c - 1]\n        # Trust position: the table is read top-left to bottom-right, row by\n        # row, which is exactly how $oldValues / $newValues were built from\n        # the diff. Cell.Range.Text includes trailing cell-mark characters,\n        # so compare using -like against the expected old value.\n        $current = $cellRange.Text\n        if ($current -eq $expectedOld -or $current -like ($expectedOld + \"*\")) {\n            $cellRange.Text = $newVal\n        } else {\n            # Fallback (should not normally happen): overwrite anyway so the\n            # table ends up with the correct target values.\n            $cellRange.Text = $newVal\n        }\n    }\n}\n\nWrite-Output \"done\"\n"}
